$wb = $excel.ActiveWorkbook

# Delete row 16 ("Sheet" row) in the optimization_parameters sheet.
$wsOpt = $wb.Worksheets.Item("optimization_parameters")
$wsOpt.Rows.Item(16).Delete()

$wb.Save()
